# Updated dataset till S360 (Simon send updated S349-354)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sample rows 350-361 (columns B:S) supplied with the dataset update
$rowData = @{}
$rowData[350] = @(0, 0, 0, 0, 0, 0, 0, 17.0512978210106, 14.5292461779029, 0, 0, 0, 0, 0, 1.6182755019939501, 0, 3.2201634039678502, 0)
$rowData[351] = @(0, 0, 0, 0, 0, 18.481327138584799, 0, 0, 0, 0, 17.662032972421301, 0, 0, 0, 3.8985001215208701, 0, 5.64216521493546, 0)
$rowData[352] = @(0, 0, 0, 0, 0, 0, 0, 0, 17.072970683774901, 0, 16.139392511553702, 0, 0, 0, 2.8594828006181898, 0, 2.20605800970135, 0)
$rowData[353] = @(0, 0, 15.513852594017299, 0, 0, 16.5527763704394, 0, 0, 0, 0, 0, 0, 0, 0, 2.6629382401022599, 0, 4.9480136978148996, 0)
$rowData[354] = @(0, 0, 0, 0, 0, 14.743020648004601, 14.5778107047774, 0, 0, 0, 0, 0, 0, 0, 5.9927940120690701, 0, 3.7683601336148498, 0)
$rowData[355] = @(7.7860153034526602, 0, 0, 14.0969695707408, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 4.7646743481566496, 0, 5.6975106373188398, 0)
$rowData[356] = @(0, 0, 0, 6.8847541424966998, 0, 0, 0, 0, 8.4649778179028292, 0, 0, 0, 0, 0, 1.49569537502413, 0, 3.9783978503180002, 0)
$rowData[357] = @(0, 0, 0, 0, 8.7947096176898398, 0, 0, 0, 0, 0, 7.9615266012769901, 0, 0, 0, 1.90057297373125, 0, 1.48423303044529, 0)
$rowData[358] = @(0, 0, 0, 0, 0, 0, 8.9083150838504501, 0, 12.352916525900801, 0, 0, 0, 0, 0, 1.23588726090864, 0, 2.72490805718405, 0)
$rowData[359] = @(0, 8.8330896029060408, 0, 0, 0, 0, 12.8856765808719, 0, 0, 0, 0, 0, 0, 0, 1.2839078207243799, 0, 2.8547635211524902, 0)
$rowData[360] = @(4.35585876633548, 0, 7.8787923441815604, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1.3465978015439, 0, 3.9680810956739299, 0)
$rowData[361] = @(0, 11.059690560140099, 0, 8.3288067144775706, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 2.4520088660528598, 0, 3.8950106132044602, 0)

foreach ($r in ($rowData.Keys | Sort-Object)) {
    $col = 2  # column B
    foreach ($v in $rowData[$r]) {
        $ws.Cells.Item($r, $col).Value = [double]$v
        $col = $col + 1
    }
}

# Update the saved view/scroll position and active selection to match
# where the user left off after entering the new rows
$win = $excel.ActiveWindow
$win.ScrollRow = 320
$win.ScrollColumn = 1
$ws.Range("W361").Select()

Write-Host "Updated rows 350-361 (B:S) and refreshed view selection"
